# "Insert new CHEMProduct regression test"
#
# For each data row (2..91), the short construct id stored in columns A and D
# is renamed from the "TEST#######" scheme to the new "TST#######" scheme
# (i.e. "TEST" -> "TST"). Column C (the long LSID, which already embeds the
# original "TEST#######" id) is left untouched, as are columns B and E:I.
#
# The sheet view's active/selected column is also moved from C to D to match
# the newly relevant column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

$colA = $ws.Range("A2:A$lastRow")
$colD = $ws.Range("D2:D$lastRow")

$valuesA = $colA.Value()
$valuesD = $colD.Value()

$rowCount = $valuesA.GetLength(0)
for ($r = 1; $r -le $rowCount; $r++) {
    $oldA = $valuesA[$r, 1]
    $newA = $oldA.Replace("TEST", "TST")
    $valuesA[$r, 1] = $newA
    $valuesD[$r, 1] = $newA
}

$colA.Value = $valuesA
$colD.Value = $valuesD

# Update the sheet view selection to match the new active column (D).
$ws.Range("D1:D1048576").Select()
